$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = ''
$ws.Range("E2").Value = ''
$ws.Range("D3").Value = ''
$ws.Range("E3").Value = 'Veg'
$ws.Range("D4").Value = ''
$ws.Range("E4").Value = 'Veg'
$ws.Range("D5").Value = ''
$ws.Range("E5").Value = 'Veg'
$ws.Range("D6").Value = ''
$ws.Range("E6").Value = ''
$ws.Range("D7").Value = ''
$ws.Range("E7").Value = ''
$ws.Range("D8").Value = 'Snack'
$ws.Range("E8").Value = 'Veg'
$ws.Range("D9").Value = 'Snack'
$ws.Range("E9").Value = ''
$ws.Range("D10").Value = 'Breakfast, Snack'
$ws.Range("E10").Value = ''
$ws.Range("D11").Value = ''
$ws.Range("E11").Value = 'Veg'
$ws.Range("D12").Value = 'Lunch, Dinner'
$ws.Range("E12").Value = ''
$ws.Range("D13").Value = ''
$ws.Range("E13").Value = 'Veg'
$ws.Range("D14").Value = ''
$ws.Range("E14").Value = 'Veg'
$ws.Range("D15").Value = 'Breakfast'
$ws.Range("E15").Value = ''
$ws.Range("D16").Value = ''
$ws.Range("E16").Value = ''
$ws.Range("D17").Value = ''
$ws.Range("E17").Value = 'Vegetarian, Veg'
$ws.Range("D18").Value = ''
$ws.Range("E18").Value = ''
$ws.Range("D19").Value = ''
$ws.Range("E19").Value = 'Veg'
$ws.Range("D20").Value = 'Snack'
$ws.Range("E20").Value = 'Jain'
$ws.Range("D21").Value = 'Snack'
$ws.Range("E21").Value = 'Veg'
$ws.Range("D22").Value = 'Breakfast'
$ws.Range("E22").Value = 'Veg'
$ws.Range("D23").Value = 'Breakfast'
$ws.Range("E23").Value = ''
$ws.Range("D24").Value = ''
$ws.Range("E24").Value = ''
$ws.Range("D25").Value = ''
$ws.Range("E25").Value = 'Vegan, Veg'
$ws.Range("D26").Value = ''
$ws.Range("E26").Value = ''
$ws.Range("D27").Value = ''
$ws.Range("E27").Value = ''
$ws.Range("D28").Value = ''
$ws.Range("E28").Value = ''
$ws.Range("D29").Value = ''
$ws.Range("E29").Value = ''
$ws.Range("D30").Value = 'Breakfast'
$ws.Range("E30").Value = ''
$ws.Range("D31").Value = ''
$ws.Range("E31").Value = ''
$ws.Range("D32").Value = 'Snack'
$ws.Range("E32").Value = 'Veg'
$ws.Range("D33").Value = 'Snack'
$ws.Range("E33").Value = 'Veg'
$ws.Range("D34").Value = 'Snack'
$ws.Range("E34").Value = 'Veg'
$ws.Range("D35").Value = ''
$ws.Range("E35").Value = 'Jain'
$ws.Range("D36").Value = 'Snack'
$ws.Range("E36").Value = 'Veg'
$ws.Range("D37").Value = 'Breakfast, Snack'
$ws.Range("E37").Value = ''
$ws.Range("D38").Value = 'Snack'
$ws.Range("E38").Value = ''
$ws.Range("D39").Value = 'Snack'
$ws.Range("E39").Value = ''
$ws.Range("D40").Value = 'Snack'
$ws.Range("E40").Value = ''
$ws.Range("D41").Value = ''
$ws.Range("E41").Value = ''
$ws.Range("D42").Value = 'Dinner'
$ws.Range("E42").Value = ''
$ws.Range("D43").Value = ''
$ws.Range("E43").Value = ''
$ws.Range("D44").Value = ''
$ws.Range("E44").Value = 'Vegan, Veg'
$ws.Range("D45").Value = ''
$ws.Range("E45").Value = ''
$ws.Range("D46").Value = 'Breakfast'
$ws.Range("E46").Value = 'Jain'
$ws.Range("D47").Value = ''
$ws.Range("E47").Value = ''
$ws.Range("D48").Value = ''
$ws.Range("E48").Value = 'Veg'
$ws.Range("D49").Value = ''
$ws.Range("E49").Value = 'Jain'
$ws.Range("D50").Value = 'Breakfast, Snack'
$ws.Range("E50").Value = ''
$ws.Range("D51").Value = ''
$ws.Range("E51").Value = ''
$ws.Range("D52").Value = ''
$ws.Range("E52").Value = 'Veg'
$ws.Range("D53").Value = 'Breakfast'
$ws.Range("E53").Value = ''
$ws.Range("D54").Value = ''
$ws.Range("E54").Value = ''
$ws.Range("D55").Value = ''
$ws.Range("E55").Value = 'Vegan, Veg'
$ws.Range("D56").Value = ''
$ws.Range("E56").Value = ''
$ws.Range("D57").Value = ''
$ws.Range("E57").Value = ''
$ws.Range("D58").Value = 'Snack'
$ws.Range("E58").Value = 'Veg'
$ws.Range("D59").Value = ''
$ws.Range("E59").Value = ''
$ws.Range("D60").Value = ''
$ws.Range("E60").Value = 'Veg'
$ws.Range("D61").Value = ''
$ws.Range("E61").Value = 'Veg'
$ws.Range("D62").Value = ''
$ws.Range("E62").Value = ''
